# Fix PER bug: team column was mis-aligned with the value column.
# Re-map each row's Team (col B) and PER value (col C) to the corrected pairing.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = "POR"
$ws.Range("C2").Value = 13.71538461538461
$ws.Range("B3").Value = "NJN"
$ws.Range("C3").Value = 12.14
$ws.Range("B4").Value = "CLE"
$ws.Range("C4").Value = 12.78461538461539
$ws.Range("B5").Value = "DAL"
$ws.Range("C5").Value = 9.916666666666666
$ws.Range("B6").Value = "MIA"
$ws.Range("C6").Value = 13.93333333333333
$ws.Range("B7").Value = "SEA"
$ws.Range("C7").Value = 13.23125
$ws.Range("B8").Value = "ATL"
$ws.Range("C8").Value = 10.1375
$ws.Range("B9").Value = "MIL"
$ws.Range("C9").Value = 9.66923076923077
$ws.Range("B10").Value = "LAC"
$ws.Range("C10").Value = 13.59285714285714
$ws.Range("B11").Value = "VAN"
$ws.Range("C11").Value = 12.32307692307692
$ws.Range("B12").Value = "DET"
$ws.Range("C12").Value = 14.29166666666667
$ws.Range("B13").Value = "WSB"
$ws.Range("C13").Value = 12.40769230769231
$ws.Range("B14").Value = "SAS"
$ws.Range("C14").Value = 14.15384615384615
$ws.Range("B15").Value = "ORL"
$ws.Range("C15").Value = 12.31428571428571
$ws.Range("B16").Value = "UTA"
$ws.Range("C16").Value = 14.74545454545454
$ws.Range("B17").Value = "HOU"
$ws.Range("C17").Value = 12.3125
$ws.Range("B18").Value = "DEN"
$ws.Range("C18").Value = 13.91
$ws.Range("B19").Value = "LAL"
$ws.Range("C19").Value = 16.25
$ws.Range("B20").Value = "GSW"
$ws.Range("C20").Value = 13.16363636363636
$ws.Range("B21").Value = "IND"
$ws.Range("C21").Value = 13.53333333333334
$ws.Range("B22").Value = "CHI"
$ws.Range("C22").Value = 14.37333333333333
$ws.Range("B23").Value = "PHI"
$ws.Range("C23").Value = 13.44285714285714
$ws.Range("B24").Value = "CHH"
$ws.Range("C24").Value = 13.22727272727273
$ws.Range("B25").Value = "BOS"
$ws.Range("C25").Value = 11.61764705882353
$ws.Range("B26").Value = "TOR"
$ws.Range("C26").Value = 11.0764705882353
$ws.Range("B27").Value = "SAC"
$ws.Range("C27").Value = 12.49285714285715
$ws.Range("B28").Value = "PHO"
$ws.Range("C28").Value = 13.31666666666667
$ws.Range("B29").Value = "NYK"
$ws.Range("C29").Value = 12.77692307692308
$ws.Range("B30").Value = "MIN"
$ws.Range("C30").Value = 13.3
